$wb = $excel.ActiveWorkbook

# --- Select C12 on the Portugal sheet (becomes the non-active selection) ---
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Activate()
$portugal.Range("C12").Select()

# Row 4's wrapped-text height is no longer needed on Portugal; rows 3 & 5 keep theirs.
$portugal.Rows.Item(4).AutoFit()

# --- Create the new "Slovakia" sheet by copying Portugal's (keeps styles / column widths / merges) ---
$portugal.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Rows 3-5 on the copied sheet inherited Portugal's taller (wrapped-text) row height;
# Slovakia's text is short, so auto-fit them back down to the sheet's default height.
$slovakia.Rows.Item(3).AutoFit()
$slovakia.Rows.Item(4).AutoFit()
$slovakia.Rows.Item(5).AutoFit()

# --- Fill in Slovakia-specific data ---
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3174"
# The value was typed directly (not copied) so it doesn't carry Portugal's cell border style.
$slovakia.Range("B4").Style = "Normal"

# --- Make Slovakia the active sheet / tab, with E15 selected ---
$slovakia.Activate()
$slovakia.Range("E15").Select()
